# Generate Report for Handback
# Updates the handoff/handback timestamp strings that get refreshed each time
# the handback status report is regenerated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
# "Latest HO Xliff Generate Date" for 63564171-3614-40b1-b3ce-7b0d8f7cbc33.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 13:08:59"

# --- zh-cn sheet -----------------------------------------------------------
# Correspond Handoff Datetime / Correspond Handback DateTime for the same file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-25 13:08:55"
$wsZhCn.Range("K2").Value = "2016-08-25 13:09:20"

# --- de-de sheet -----------------------------------------------------------
# Correspond Handoff Datetime / Correspond Handback DateTime for the same file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-25 13:08:59"
$wsDeDe.Range("K2").Value = "2016-08-25 13:09:27"
